$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I4").Value = 0.016
